$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$msgKids = @"
🔥 **DEMO DEALS OF THE DAY** 🔥
💥 Price Crash Store
⚡️ Up to 5% off
👉🏻 amzaff.in/l3swo0g
🌟 Kid's Carnival
📣 Sale live till 28th July
👉🏻 amzaff.in/jrtPYsT
🧸 Toy's Fiesta
⚡️ Up to 70% off
👉🏻 amzaff.in/pXpqAHe
🪴 Solar Garden & Outdoors
⚡️ Up to 70% off
👉🏻 amzaff.in/092HVM2
🛍 End Of Season Sale
⚡️ Up to 70% off
👉🏻 amzaff.in/ukuS1wj
🏠 Home Shopping Spree
📣 Sale live from 26th-30th July
👉🏻 amzaff.in/dvgN1JH
"@

$msgEssentials = @"
💧 **DEMO Essentials for Skin, Hair & Fragrance** 💧
💥 Maximise earnings with Beauty commissions – Now increased to 10%
🌿 Mamaearth Rice Oil-Free Face Moisturizer for Oily Skin
⚡️ 80g @25% + 5% Off – ₹22
👉🏻 amzaff.in/kyKGkVq
🧼 Cetaphil Brightness Reveal Creamy Cleanser
⚡️ 100g @25% Off – ₹599
👉🏻 amzaff.in/2QkCAT6
💆‍♀️ Herbal Essences bio:renew Argan Oil of Morocco Shampoo
⚡️ 400ml @60% Off – ₹260
👉🏻 amzaff.in/WPdHenG
🧴 Be Bodywise 6% AHA BHA Underarm Roll On Deodorant
⚡️ 50ml | Alcohol & Aluminum Free – ₹399
👉🏻 amzaff.in/3YBXpxC
🕺 Park Avenue Voyage Signature Collection Perfume for Men
⚡️ 120ml @57% Off – ₹171
👉🏻 amzaff.in/IQsEdXX
"@

$msgLaptops = @"
💻 **DEMO Level Gaming Laptops**
⚡️ Up to 45% off
👉🏻 amzaff.in/FeVABNi
🔥 Gaming Laptops Under ₹80,000 🔥
💻 ASUS TUF A15 GAMING
⚡️ ~~₹83,990~~ | ₹63,990
⚡️ Effective price ₹56,490
👉🏻 amzaff.in/HNJJ4b3
💻 HP VICTUS GAMING
⚡️ ~~₹99,382~~ | ₹82,990
⚡️ Effective price ₹77,490
👉🏻 amzaff.in/rIYTx8U
💻 LENOVO LOQ GAMING
⚡️ ~~₹1,12,990~~ | ₹87,190
⚡️ Effective price ₹77,190
👉🏻 amzaff.in/Q2oXx7g
💻 ACER NITRO V GAMING
⚡️ ~~₹89,999~~ | ₹70,990
⚡️ Effective price ₹67,490
👉🏻 amzaff.in/Es6mtU0
💻 DELL G15 GAMING
⚡️ ~~₹1,06,331~~ | ₹77,490
⚡️ Effective price ₹68,990
👉🏻 amzaff.in/MMEYXBc
💻 ASUS CREATOR SERIES
⚡️ ~~₹85,990~~ | ₹72,990
⚡️ Effective price ₹66,490
👉🏻 amzaff.in/R4IV7C8
"@

$rows = @(
    @{ Row = 5;  A = 1; C = "2025-07-31 14:30:00"; D = "Kid's Carnival";    E = $msgKids },
    @{ Row = 6;  A = 2; C = "2025-07-31 16:00:00"; D = "Daily Essentials";   E = $msgEssentials },
    @{ Row = 7;  A = 3; C = "2025-07-31 17:30:00"; D = "Laptops";            E = $msgLaptops },
    @{ Row = 8;  A = 1; C = "2025-07-31 14:30:00"; D = "Kid's Carnival";    E = $msgKids },
    @{ Row = 9;  A = 2; C = "2025-07-31 16:00:00"; D = "Daily Essentials";   E = $msgEssentials },
    @{ Row = 10; A = 3; C = "2025-07-31 17:30:00"; D = "Laptops";            E = $msgLaptops }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = "✅ Scheduled"
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
}

